$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'78.925.14"
$ws.Range("E2").Value = "  +3.23%  "
$ws.Range("D3").Value = "'3.193.24"
$ws.Range("E3").Value = "  +5.42%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'205.53"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").Value = "'632.18"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +10.23%  "
$ws.Range("E9").Value = "  +5.30%  "
$ws.Range("D10").Value = "'3.193.10"
$ws.Range("E10").Value = "  +5.43%  "
$ws.Range("D11").Value = "'0.579"
$ws.Range("E11").Value = "  +32.94%  "
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").Value = "'5.42"
$ws.Range("E13").Value = "  +6.25%  "
$ws.Range("D14").Value = "'3.775.97"
$ws.Range("E14").Value = "  +5.24%  "
$ws.Range("D15").Value = "'0.0000223"
$ws.Range("E15").Value = "  +17.67%  "
$ws.Range("D16").Value = "'31.42"
$ws.Range("E16").Value = "  +7.83%  "
$ws.Range("D17").Value = "'78.921.50"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").Value = "'3.190.38"
$ws.Range("E18").Value = "  +5.06%  "
$ws.Range("D19").Value = "'14.45"
$ws.Range("E19").Value = "  +7.42%  "
$ws.Range("D20").Value = "'9.29"
$ws.Range("E20").Value = "  +2.60%  "
$ws.Range("D21").Value = "'427.59"
$ws.Range("E21").Value = "  +14.61%  "
$ws.Range("E22").Value = "  +26.00%  "
$ws.Range("E23").Value = "  +12.81%  "
$ws.Range("D24").Value = "'6.85"
$ws.Range("E24").Value = "  +5.72%  "
$ws.Range("D25").Value = "'3.349.57"
$ws.Range("E25").Value = "  +5.11%  "
$ws.Range("E26").Value = "  +8.43%  "
$ws.Range("D27").Value = "'11.04"
$ws.Range("E27").Value = "  +11.65%  "
$ws.Range("D28").Value = "'75.64"
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'8.81"
$ws.Range("E32").Value = "  +6.23%  "
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("D34").Value = "'513.05"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").Value = "'1.99"
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("E36").Value = "  +22.36%  "
$ws.Range("D37").Value = "'22.90"
$ws.Range("E37").Value = "  +10.55%  "
$ws.Range("E38").Value = "  +19.85%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "'0.398"
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("D41").Value = "'163.96"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'19.98"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "'192.31"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  +6.65%  "
$ws.Range("E46").Value = "  +12.56%  "
$ws.Range("E47").Value = "  +7.06%  "
$ws.Range("D48").Value = "'1.29"
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("D49").Value = "'42.68"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "'2.48"
$ws.Range("E50").Value = "  +6.24%  "
$ws.Range("D51").Value = "'0.621"
$ws.Range("E51").Value = "  +3.12%  "
